# Atualização dos dados e melhorias no codigo
#
# - Insert a new header row at the top of the sheet (pushes all existing
#   data rows down by one).
# - Populate the new header row with the column titles.
# - Format the header row: bold font, thin border all around the cells,
#   text centered horizontally and aligned to the top vertically.
# - Fix a typo in the municipality name ("ilhabela" -> "ilha bela"),
#   which after the row insert now lives on row 53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the existing row 1, shifting all data down.
$ws.Rows.Item(1).Insert()

# New header labels.
$ws.Range("A1").Value = "MUNICIPIO"
$ws.Range("B1").Value = "Nº DE CASOS"
$ws.Range("C1").Value = "ÓBITOS"

# Header formatting: thin border box, bold text, centered / top-aligned.
$headerRange = $ws.Range("A1:C1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Correct the municipality name typo (now shifted to row 53).
$ws.Range("A53").Value = "ilha bela"
